$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 20.268639
$ws.Range("H2").Value = 60.80591700000001
$ws.Range("I2").Value = 0.09537690375401175
$ws.Range("J2").Value = 0.09537690375401174
$ws.Range("M2").Value = 0.01989833333333333
$ws.Range("N2").Value = 0.059695
$ws.Range("O2").Value = 0.5455534129646046
$ws.Range("P2").Value = 0.5455534129646046
$ws.Range("Q2").Value = 0.4033121350350001
$ws.Range("R2").Value = 3.629809215315
$ws.Range("S2").Value = 0.05203319536099772
$ws.Range("T2").Value = 0.05203319536099772

$ws.Range("G3").Value = 20.268639
$ws.Range("H3").Value = 60.80591700000001
$ws.Range("I3").Value = 0.09537690375401175
$ws.Range("J3").Value = 0.09537690375401174
$ws.Range("O3").Value = 0.1997148627777118
$ws.Range("P3").Value = 0.1997148627777118
$ws.Range("Q3").Value = 0.147643522689
$ws.Range("R3").Value = 1.328791704201
$ws.Range("S3").Value = 0.01904818524539548
$ws.Range("T3").Value = 0.01904818524539548

$ws.Range("G4").Value = 20.268639
$ws.Range("H4").Value = 60.80591700000001
$ws.Range("I4").Value = 0.09537690375401175
$ws.Range("J4").Value = 0.09537690375401174
$ws.Range("M4").Value = 0.009290999999999999
$ws.Range("N4").Value = 0.027873
$ws.Range("O4").Value = 0.2547317242576836
$ws.Range("P4").Value = 0.2547317242576836
$ws.Range("Q4").Value = 0.188315924949
$ws.Range("R4").Value = 1.694843324541
$ws.Range("S4").Value = 0.02429552314761855
$ws.Range("T4").Value = 0.02429552314761855

$ws.Range("I5").Value = 0.6011031624655011
$ws.Range("J5").Value = 0.601103162465501
$ws.Range("M5").Value = 0.01989833333333333
$ws.Range("N5").Value = 0.059695
$ws.Range("O5").Value = 0.5455534129646046
$ws.Range("P5").Value = 0.5455534129646046
$ws.Range("Q5").Value = 2.541833402932778
$ws.Range("R5").Value = 22.876500626395
$ws.Range("S5").Value = 0.3279338818268713
$ws.Range("T5").Value = 0.3279338818268713

$ws.Range("I6").Value = 0.6011031624655011
$ws.Range("J6").Value = 0.601103162465501
$ws.Range("O6").Value = 0.1997148627777118
$ws.Range("P6").Value = 0.1997148627777118
$ws.Range("S6").Value = 0.1200492356070461
$ws.Range("T6").Value = 0.1200492356070461

$ws.Range("I7").Value = 0.6011031624655011
$ws.Range("J7").Value = 0.601103162465501
$ws.Range("M7").Value = 0.009290999999999999
$ws.Range("N7").Value = 0.027873
$ws.Range("O7").Value = 0.2547317242576836
$ws.Range("P7").Value = 0.2547317242576836
$ws.Range("Q7").Value = 1.186841819917
$ws.Range("R7").Value = 10.681576379253
$ws.Range("S7").Value = 0.1531200450315836
$ws.Range("T7").Value = 0.1531200450315836

$ws.Range("G8").Value = 42.02733833333333
$ws.Range("H8").Value = 126.082015
$ws.Range("I8").Value = 0.1977654939365007
$ws.Range("J8").Value = 0.1977654939365007
$ws.Range("M8").Value = 0.01989833333333333
$ws.Range("N8").Value = 0.059695
$ws.Range("O8").Value = 0.5455534129646046
$ws.Range("P8").Value = 0.5455534129646046
$ws.Range("Q8").Value = 0.8362739872694444
$ws.Range("R8").Value = 7.526465885424999
$ws.Range("S8").Value = 0.1078916401836888
$ws.Range("T8").Value = 0.1078916401836888

$ws.Range("G9").Value = 42.02733833333333
$ws.Range("H9").Value = 126.082015
$ws.Range("I9").Value = 0.1977654939365007
$ws.Range("J9").Value = 0.1977654939365007
$ws.Range("O9").Value = 0.1997148627777118
$ws.Range("P9").Value = 0.1997148627777118
$ws.Range("Q9").Value = 0.3061411415327777
$ws.Range("R9").Value = 2.755270273795
$ws.Range("S9").Value = 0.03949670848369462
$ws.Range("T9").Value = 0.03949670848369462

$ws.Range("G10").Value = 42.02733833333333
$ws.Range("H10").Value = 126.082015
$ws.Range("I10").Value = 0.1977654939365007
$ws.Range("J10").Value = 0.1977654939365007
$ws.Range("M10").Value = 0.009290999999999999
$ws.Range("N10").Value = 0.027873
$ws.Range("O10").Value = 0.2547317242576836
$ws.Range("P10").Value = 0.2547317242576836
$ws.Range("Q10").Value = 0.3904760004549999
$ws.Range("R10").Value = 3.514284004094999
$ws.Range("S10").Value = 0.05037714526911728
$ws.Range("T10").Value = 0.0503771452691173

$ws.Range("G11").Value = 22.47397933333333
$ws.Range("H11").Value = 67.421938
$ws.Range("I11").Value = 0.1057544398439867
$ws.Range("J11").Value = 0.1057544398439867
$ws.Range("M11").Value = 0.01989833333333333
$ws.Range("N11").Value = 0.059695
$ws.Range("O11").Value = 0.5455534129646046
$ws.Range("P11").Value = 0.5455534129646046
$ws.Range("Q11").Value = 0.4471947321011111
$ws.Range("R11").Value = 4.024752588909999
$ws.Range("S11").Value = 0.05769469559304689
$ws.Range("T11").Value = 0.05769469559304689

$ws.Range("G12").Value = 22.47397933333333
$ws.Range("H12").Value = 67.421938
$ws.Range("I12").Value = 0.1057544398439867
$ws.Range("J12").Value = 0.1057544398439867
$ws.Range("O12").Value = 0.1997148627777118
$ws.Range("P12").Value = 0.1997148627777118
$ws.Range("Q12").Value = 0.1637079567904444
$ws.Range("R12").Value = 1.473371611114
$ws.Range("S12").Value = 0.02112073344157557
$ws.Range("T12").Value = 0.02112073344157557

$ws.Range("G13").Value = 22.47397933333333
$ws.Range("H13").Value = 67.421938
$ws.Range("I13").Value = 0.1057544398439867
$ws.Range("J13").Value = 0.1057544398439867
$ws.Range("M13").Value = 0.009290999999999999
$ws.Range("N13").Value = 0.027873
$ws.Range("O13").Value = 0.2547317242576836
$ws.Range("P13").Value = 0.2547317242576836
$ws.Range("Q13").Value = 0.208805741986
$ws.Range("R13").Value = 1.879251677874
$ws.Range("S13").Value = 0.0269390108093642
$ws.Range("T13").Value = 0.0269390108093642
